$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new formula to C2 (mirrors C1's formula, but referencing A2 instead of A1)
$ws.Range("C2").Formula = "=selectif(2,A1=1,A2+B1,)"

# Move/update the active selection to C2
$ws.Range("C2").Select()
